$d = $word.ActiveDocument

$pairs = @(
    @("62×48=", "24×79="),
    @("87×34=", "70×11="),
    @("66×20=", "26×99="),
    @("76×38=", "11×47="),
    @("51×62=", "89×56="),
    @("33×26=", "58×68="),
    @("13×83=", "79×89="),
    @("38×90=", "99×40="),
    @("53×43=", "67×32="),
    @("37×18=", "33×67="),
    @("41×83=", "67×82="),
    @("55×82=", "45×66="),
    @("98×65=", "95×94="),
    @("47×84=", "86×25="),
    @("91×30=", "71×88="),
    @("68×95=", "59×39="),
    @("17×88=", "92×94="),
    @("72×40=", "35×63="),
    @("91×66=", "45×62="),
    @("64×87=", "43×14="),
    @("51×73=", "50×39="),
    @("11×72=", "53×27="),
    @("44×68=", "40×19="),
    @("68×31=", "54×76="),
    @("46×51=", "45×41=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
